# fall 24 week 10 inputs and lineup message improvements
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows of matchup data to append (columns A:D), starting at row 1917
$newData = @(
    @(5,8,6,12),
    @(4,16,3,4),
    @(6,4,9,16),
    @(4,5,3,15),
    @(5,7,4,13),
    @(3,16,4,4),
    @(4,3,5,17),
    @(7,13,4,7),
    @(3,6,2,14),
    @(5,15,7,5),
    @(5,5,4,15),
    @(5,14,7,6),
    @(3,5,2,15),
    @(5,4,4,16),
    @(5,15,4,5),
    @(4,2,3,18),
    @(3,13,1,7),
    @(4,18,3,2),
    @(5,16,4,4),
    @(6,8,5,12),
    @(9,5,3,15),
    @(3,8,5,12),
    @(3,15,4,5),
    @(4,6,5,14),
    @(4,7,5,13),
    @(7,7,5,13),
    @(5,6,4,14),
    @(4,5,3,15),
    @(6,12,7,8),
    @(4,6,2,14)
)

$startRow = 1917
$rowCount = $newData.Count
$endRow = $startRow + $rowCount - 1

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $startRow + $i
    $rowVals = $newData[$i]
    $ws.Cells.Item($r, 1).Value = $rowVals[0]
    $ws.Cells.Item($r, 2).Value = $rowVals[1]
    $ws.Cells.Item($r, 3).Value = $rowVals[2]
    $ws.Cells.Item($r, 4).Value = $rowVals[3]
}

# Update the view so the new last row is visible/selected, mirroring the
# resulting sheetView state after the data was appended in Excel.
$excel.ActiveWindow.ScrollRow = $startRow + 14
$selCell = "A" + ($endRow + 1)
$ws.Range($selCell).Select()
